$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.658.45"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "'2.605.69"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'572.45"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "'154.73"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -5.81%  "
$ws.Range("E9").Value = "  -6.36%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'28.07"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").Value = "'3.072.17"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("E15").Value = "  -8.28%  "
$ws.Range("D16").Value = "'63.404.16"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").Value = "'2.548.84"
$ws.Range("E17").Value = "  -4.26%  "
$ws.Range("D18").Value = "'11.98"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("D19").Value = "'7.54"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "'4.53"
$ws.Range("E20").Value = "  -5.08%  "
$ws.Range("D21").Value = "'341.68"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D23").Value = "'66.89"
$ws.Range("E23").Value = "  -4.04%  "
$ws.Range("D24").Value = "'1.77"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("D26").Value = "'587.25"
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").Value = "'1.73"
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("D34").Value = "'6.50"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "'5.36"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").Value = "'0.404"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'19.59"
$ws.Range("E38").Value = "  -4.62%  "
$ws.Range("D39").Value = "'153.81"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'41.34"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").Value = "'2.43"
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("D44").Value = "'155.33"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").Value = "'23.05"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("D47").Value = "'0.0585"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").Value = "'18.75"
$ws.Range("E51").Value = "  -5.26%  "
